$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-10 Monday" "2024-06-11 Tuesday"

Replace-Text "26×25=" "96×74="
Replace-Text "56×29=" "48×78="
Replace-Text "65×81=" "64×57="
Replace-Text "18×42=" "28×15="
Replace-Text "61×20=" "90×85="

Replace-Text "43×42=" "46×39="
Replace-Text "58×62=" "74×19="
Replace-Text "97×17=" "34×55="
Replace-Text "90×94=" "36×27="
Replace-Text "78×26=" "77×90="

Replace-Text "94×17=" "48×99="
Replace-Text "36×99=" "31×29="
Replace-Text "62×53=" "65×59="
Replace-Text "66×58=" "60×70="
Replace-Text "43×26=" "61×77="

Replace-Text "27×89=" "95×69="
Replace-Text "47×67=" "77×32="
Replace-Text "74×87=" "60×86="
Replace-Text "50×34=" "74×13="
Replace-Text "16×68=" "60×41="

Replace-Text "24×62=" "29×32="
Replace-Text "86×64=" "50×55="
Replace-Text "74×70=" "77×97="
Replace-Text "80×34=" "36×20="
Replace-Text "25×78=" "19×90="
